$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet, positioned right before the
#    "总计" (Total) summary sheet, by copying the "2021-Q4" sheet as a
#    template (same column layout / header / styles) and then replacing
#    its data rows with the 2022-Q1 fund-holding figures.
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Sheet references above resolve positionally, so re-fetch "总计" by name
# now that the sheet collection has shifted.
$totalSheet = $wb.Worksheets.Item("总计")

# The template ("2021-Q4") has 8 data rows (rows 2-9); 2022-Q1 only has
# 6, so drop the two extra rows before writing new content.
$newSheet.Rows("8:9").Delete()

$fundRows = @(
    ,@("506001", "万家科创板 2 年定期开放混合型证券投资基金", "12.84", "98.14", "3.59", "0.4610", 6)
    ,@("161605", "融通蓝筹成长混合", "4.82", "71.70", "5.32", "0.2564", 3)
    ,@("000717", "融通转型三动力灵活配置混合A", "3.83", "94.89", "6.29", "0.2409", 3)
    ,@("009828", "融通转型三动力灵活配置混合C", "0.59", "94.89", "6.29", "0.0371", 3)
    ,@("004917", "中银证券祥瑞混合A", "0.10", "79.01", "2.53", "0.0025", 5)
    ,@("004918", "中银证券祥瑞混合C", "0.07", "79.01", "2.53", "0.0018", 5)
)

# Text-valued columns (B..G) must stay text even though some look numeric
# (fund codes with leading zeros, decimal figures stored as strings).
$newSheet.Range("B2:G7").NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = 2 + $i
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new first data row for 2022-Q1
#    (above the existing 2021-Q4 row) and renumber the index column.
# ------------------------------------------------------------------
$totalSheet.Rows("2:2").Insert()

# The inserted row picks up the header row's bold/centred formatting;
# reset it back to normal before writing the new values.
$totalSheet.Range("B2:D2").Style = "Normal"

# Column A's index cells use a distinct style (bold + thin border); copy
# it verbatim from the row below rather than re-deriving it by hand.
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 1

# Renumber the index column (A) for the remaining rows (they each shift
# down by one position, 0-based).
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally active sheet/selection (the sheet Copy() above
# leaves "2022-Q1" active).
$wb.Worksheets.Item("2020-Q4").Activate()
